$d = $word.ActiveDocument
$d.Content.Find.Execute("GeomAlg2019Oct", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "GeomAlg2021Jan", 2)
